# Add five new attribute rows (40-44) to the features sheet, matching the
# existing "directional damage" style rows already present (e.g. row 38/39):
# column A left blank, B = attribute name, C = options, D = reasoning,
# E = uncertainty, F = hazard relevance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 40; A = ""; B = "Wall Damage % (Front)"; C = "0%, <25%, 25-50%, 50-75%, >75%"; D = 'Analyze "front" images. Estimate % wall failure. If obscured/missing, use [NEEDS_RESEARCH].'; E = "Low"; F = "All" },
    @{ Row = 41; A = ""; B = "Wall Damage % (Rear)";  C = "0%, <25%, 25-50%, 50-75%, >75%"; D = 'Analyze "rear" images. Estimate % wall failure. If obscured/missing, use [NEEDS_RESEARCH].';  E = "Low"; F = "All" },
    @{ Row = 42; A = ""; B = "Wall Damage % (Left)";  C = "0%, <25%, 25-50%, 50-75%, >75%"; D = 'Analyze "left" images. Estimate % wall failure. If obscured/missing, use [NEEDS_RESEARCH].';  E = "Low"; F = "All" },
    @{ Row = 43; A = ""; B = "Wall Damage % (Right)"; C = "0%, <25%, 25-50%, 50-75%, >75%"; D = 'Analyze "right" images. Estimate % wall failure. If obscured/missing, use [NEEDS_RESEARCH].'; E = "Low"; F = "All" },
    @{ Row = 44; A = ""; B = "Flood Duration (Hours)"; C = "Numeric"; D = "Check specific reports or interviews. Very hard to tell visually. Likely [NEEDS_RESEARCH]."; E = "High"; F = "Flood" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
